$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = "Normal"
}

Set-TextValue "D2" "272.63"
Set-TextValue "E2" "-0.20%"
Set-TextValue "D3" "26.70"
Set-TextValue "E3" "-0.67%"
Set-TextValue "D4" "4.900"
Set-TextValue "E4" "3.88%"
Set-TextValue "D5" "0.06322"
Set-TextValue "E5" "2.61%"
Set-TextValue "E6" "2.35%"
Set-TextValue "D7" "3.356"
Set-TextValue "E7" "5.29%"
Set-TextValue "D8" "1.330"
Set-TextValue "E8" "45.90%"
Set-TextValue "D9" "0.8864"
Set-TextValue "E9" "2.91%"
Set-TextValue "E10" "1.57%"
Set-TextValue "D11" "0.05163"
Set-TextValue "E11" "-2.98%"
Set-TextValue "D12" "0.07391"
Set-TextValue "E12" "3.18%"
Set-TextValue "D13" "0.03128"
Set-TextValue "E13" "-1.83%"
Set-TextValue "D14" "0.09036"
Set-TextValue "E14" "-0.10%"
Set-TextValue "D15" "0.001574"
Set-TextValue "E15" "1.69%"
Set-TextValue "D16" "0.0006304"
Set-TextValue "E16" "3.48%"
Set-TextValue "E17" "1.22%"
Set-TextValue "D18" "3.461"
Set-TextValue "E18" "-0.24%"
Set-TextValue "D19" "2.271"
Set-TextValue "E19" "0.29%"
Set-TextValue "E20" "2.61%"
Set-TextValue "D21" "0.1316"
Set-TextValue "E21" "0.54%"
Set-TextValue "D22" "3.914"
Set-TextValue "E22" "1.75%"
Set-TextValue "D23" "0.04362"
Set-TextValue "E23" "2.47%"
Set-TextValue "D24" "0.001178"
Set-TextValue "E24" "0.18%"
Set-TextValue "D25" "0.003675"
Set-TextValue "E25" "-12.25%"
Set-TextValue "D26" "0.0001201"
Set-TextValue "E26" "0.16%"
Set-TextValue "D27" "0.0001697"
Set-TextValue "E27" "1.31%"
Set-TextValue "D40" "0.04029"
Set-TextValue "E40" "1.37%"
Set-TextValue "D41" "0.006616"
Set-TextValue "E41" "6.22%"
Set-TextValue "D42" "0.1162"
Set-TextValue "D43" "0.002102"
Set-TextValue "E43" "-3.07%"
Set-TextValue "D44" "0.01225"
Set-TextValue "E44" "-4.59%"
Set-TextValue "D45" "0.00005315"
Set-TextValue "E45" "3.62%"
Set-TextValue "D46" "2.362"
Set-TextValue "E46" "163.54%"
Set-TextValue "D47" "0.02600"
Set-TextValue "E47" "-13.02%"
